$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update name / email data (row 2 & 3) ---
$ws.Range("B2").Value = "Arun Nair"
$ws.Range("C2").Value = "appuarunnair@gmail.com"
$ws.Range("D2").Value = 9820055038

$ws.Range("B3").Value = "Aditya Kurhade"
$ws.Range("C3").Value = "adikrhd@gmail.com"
$ws.Range("D3").Value = 9664240554

# --- Add new UID column (E) ---
$ws.Range("E1").Value = "UID"
$ws.Range("E2").Value = "CSH14030"
$ws.Range("E3").Value = "CSH14038"

# --- Re-create the two hyperlinks that still apply (emails changed) and drop the 3rd (row 4 removed) ---
$ws.Hyperlinks.Delete()
$ws.Hyperlinks.Add($ws.Range("C2"), "mailto:appuarunnair@gmail.com")
$ws.Hyperlinks.Add($ws.Range("C3"), "mailto:adikrhd@gmail.com")
$ws.Range("C2").Style = "Hyperlink"
$ws.Range("C3").Style = "Hyperlink"

# --- Header formatting: center align + green fill across A1:E1 ---
$ws.Range("A1:E1").HorizontalAlignment = -4108
$ws.Range("A1:E1").Interior.Color = 5296274

# --- Remove the third student row (row 4), leaving only the hyperlink-styled, empty C4 ---
$ws.Range("A4:E4").ClearContents()
$ws.Range("C4").Style = "Hyperlink"

# --- Dimension / selection bookkeeping to match the edited workbook ---
$ws.Range("A4").Select()
